$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 5885.7144
$ws.Range("I40").Value = 5885.7144
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5885.7144
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5710.7144
$ws.Range("N40").ClearContents()
# Row 49
$ws.Range("H49").Value = 385
$ws.Range("I49").Value = 555
$ws.Range("J49").Value = 300
$ws.Range("K49").Value = 1665
$ws.Range("L49").Value = 900
$ws.Range("M49").Value = -1529
$ws.Range("N49").Value = -1172
# Row 58
$ws.Range("H58").Value = 507.8889
$ws.Range("I58").Value = 338.7143
$ws.Range("J58").Value = 1100
$ws.Range("K58").Value = 1016.1429
$ws.Range("L58").Value = 3300
$ws.Range("M58").Value = -866.1428999999999
$ws.Range("N58").Value = -3600
# Row 62
$ws.Range("H62").Value = 1983.3334
$ws.Range("I62").Value = 983.3333
$ws.Range("J62").Value = 2983.3333
$ws.Range("K62").Value = 983.3333
$ws.Range("L62").Value = 2983.3333
$ws.Range("M62").Value = -359.3333
$ws.Range("N62").Value = -4231.3333
# Row 65
$ws.Range("H65").Value = 1983.3334
$ws.Range("I65").Value = 983.3333
$ws.Range("J65").Value = 2983.3333
$ws.Range("K65").Value = 4916.6665
$ws.Range("L65").Value = 14916.6665
$ws.Range("M65").Value = -1796.6665
$ws.Range("N65").Value = -21156.6665
# Row 86
$ws.Range("H86").Value = 3356.3667
$ws.Range("J86").Value = 2881.7856
$ws.Range("L86").Value = 2881.7856
$ws.Range("N86").Value = -5127.7856
# Row 89
$ws.Range("H89").Value = 3356.3667
$ws.Range("J89").Value = 2881.7856
$ws.Range("L89").Value = 14408.928
$ws.Range("N89").Value = -25640.928
# Row 106
$ws.Range("H106").Value = 2945.8108
$ws.Range("I106").Value = 2713.5715
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 2713.5715
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -2082.5715
$ws.Range("N106").Value = -4262

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4002.281
$ws.Range("I32").Value = 2906.1948
$ws.Range("J32").Value = 11035.5
$ws.Range("K32").Value = 2906.1948
$ws.Range("L32").Value = 11035.5
$ws.Range("M32").Value = -2619.1948
$ws.Range("N32").Value = -11609.5
# Row 74
$ws.Range("H74").Value = 1888.2609
$ws.Range("I74").Value = 824.8570999999999
$ws.Range("J74").Value = 3542.4443
$ws.Range("K74").Value = 824.8570999999999
$ws.Range("L74").Value = 3542.4443
$ws.Range("M74").Value = 49.14290000000005
$ws.Range("N74").Value = -5290.4443
# Row 77
$ws.Range("H77").Value = 1888.2609
$ws.Range("I77").Value = 824.8570999999999
$ws.Range("J77").Value = 3542.4443
$ws.Range("K77").Value = 4124.2855
$ws.Range("L77").Value = 17712.2215
$ws.Range("M77").Value = 243.7145
$ws.Range("N77").Value = -26448.2215

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 12150
$ws.Range("I2").Value = 966.6667
$ws.Range("J2").Value = 23333.334
$ws.Range("K2").Value = 966.6667
$ws.Range("L2").Value = 23333.334
$ws.Range("M2").Value = -853.6667
$ws.Range("N2").Value = -23559.334
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 3361
$ws.Range("J34").Value = 4000
$ws.Range("L34").Value = 12000
$ws.Range("N34").Value = -12168
# Row 55
$ws.Range("H55").Value = 2041.05
$ws.Range("I55").Value = 868.44446
$ws.Range("J55").Value = 3000.4546
$ws.Range("K55").Value = 2605.33338
$ws.Range("L55").Value = 9001.363799999999
$ws.Range("M55").Value = -2428.33338
$ws.Range("N55").Value = -9355.363799999999
# Row 64
$ws.Range("H64").Value = 3920.6667
$ws.Range("I64").Value = 762
$ws.Range("J64").Value = 5500
$ws.Range("K64").Value = 2286
$ws.Range("L64").Value = 16500
$ws.Range("M64").Value = -2016
$ws.Range("N64").Value = -17040
# Row 67
$ws.Range("H67").Value = 3920.6667
$ws.Range("I67").Value = 762
$ws.Range("J67").Value = 5500
$ws.Range("K67").Value = 2286
$ws.Range("L67").Value = 16500
$ws.Range("M67").Value = -1350
$ws.Range("N67").Value = -18372
# Row 116
$ws.Range("H116").Value = 6735.0835
$ws.Range("I116").Value = 5997.8
$ws.Range("K116").Value = 17993.4
$ws.Range("M116").Value = -14551.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 7554.3335
$ws.Range("I5").Value = 2752
$ws.Range("J5").Value = 8926.429
$ws.Range("K5").Value = 2752
$ws.Range("L5").Value = 8926.429
$ws.Range("M5").Value = -2640
$ws.Range("N5").Value = -9150.429
# Row 70
$ws.Range("H70").Value = 4226.304
$ws.Range("I70").Value = 4269.6924
$ws.Range("J70").Value = 4169.9
$ws.Range("K70").Value = 4269.6924
$ws.Range("L70").Value = 4169.9
$ws.Range("M70").Value = -3999.6924
$ws.Range("N70").Value = -4709.9
# Row 73
$ws.Range("H73").Value = 4226.304
$ws.Range("I73").Value = 4269.6924
$ws.Range("J73").Value = 4169.9
$ws.Range("K73").Value = 4269.6924
$ws.Range("L73").Value = 4169.9
$ws.Range("M73").Value = -3333.6924
$ws.Range("N73").Value = -6041.9
# Row 80
$ws.Range("H80").Value = 2237.5
$ws.Range("I80").Value = 2170
$ws.Range("J80").Value = 2350
$ws.Range("K80").Value = 2170
$ws.Range("L80").Value = 2350
$ws.Range("M80").Value = -1172
$ws.Range("N80").Value = -4346
# Row 83
$ws.Range("H83").Value = 2237.5
$ws.Range("I83").Value = 2170
$ws.Range("J83").Value = 2350
$ws.Range("K83").Value = 10850
$ws.Range("L83").Value = 11750
$ws.Range("M83").Value = -5858
$ws.Range("N83").Value = -21734

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 35
$ws.Range("H35").Value = 10425
$ws.Range("I35").Value = 541.4286
$ws.Range("J35").Value = 45017.5
$ws.Range("K35").Value = 541.4286
$ws.Range("L35").Value = 45017.5
$ws.Range("M35").Value = -205.4286
$ws.Range("N35").Value = -45689.5
# Row 46
$ws.Range("H46").Value = 666
$ws.Range("I46").Value = 699
$ws.Range("K46").Value = 699
$ws.Range("M46").Value = -511
# Row 122
$ws.Range("H122").Value = 6884
$ws.Range("I122").Value = 11250
$ws.Range("J122").Value = 4701
$ws.Range("K122").Value = 33750
$ws.Range("L122").Value = 14103
$ws.Range("M122").Value = -31300
$ws.Range("N122").Value = -19003

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 13
$ws.Range("H13").Value = 20003
$ws.Range("J13").Value = 20003
$ws.Range("L13").Value = 20003
$ws.Range("N13").Value = -20283
